$d = $word.ActiveDocument

$replacements = @(
    @{old="792÷9=88, 0"; new="631÷9=70, 1"},
    @{old="675÷6=112, 3"; new="915÷7=130, 5"},
    @{old="100÷3=33, 1"; new="367÷3=122, 1"},
    @{old="865÷9=96, 1"; new="133÷3=44, 1"},
    @{old="192÷4=48, 0"; new="786÷8=98, 2"},
    @{old="988÷9=109, 7"; new="866÷4=216, 2"},
    @{old="486÷3=162, 0"; new="726÷2=363, 0"},
    @{old="357÷2=178, 1"; new="447÷2=223, 1"},
    @{old="987÷9=109, 6"; new="660÷5=132, 0"},
    @{old="373÷6=62, 1"; new="592÷8=74, 0"},
    @{old="448÷9=49, 7"; new="338÷9=37, 5"},
    @{old="425÷2=212, 1"; new="783÷8=97, 7"},
    @{old="805÷3=268, 1"; new="491÷8=61, 3"},
    @{old="318÷6=53, 0"; new="868÷2=434, 0"},
    @{old="787÷3=262, 1"; new="709÷3=236, 1"},
    @{old="932÷5=186, 2"; new="370÷2=185, 0"},
    @{old="315÷9=35, 0"; new="417÷3=139, 0"},
    @{old="402÷5=80, 2"; new="893÷9=99, 2"},
    @{old="295÷3=98, 1"; new="577÷8=72, 1"},
    @{old="373÷7=53, 2"; new="695÷3=231, 2"},
    @{old="537÷8=67, 1"; new="646÷8=80, 6"},
    @{old="603÷5=120, 3"; new="746÷2=373, 0"},
    @{old="502÷2=251, 0"; new="801÷2=400, 1"},
    @{old="426÷5=85, 1"; new="319÷6=53, 1"},
    @{old="501÷5=100, 1"; new="636÷8=79, 4"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
